$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = 'Brusoni et al_2001_ASQ_Knowledge Specialization Organizational Coupling and the Boundaries of the Firm.pdf'
$ws.Range("B31").Value = 'F1_P10_Brusoni et al_2001_ASQ_Knowledge Specialization Organizational Coupling and the Boundaries of the Firm.png'
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 2001
$ws.Range("E31").Value = 11
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 'final_figures/2001/F1_P10_Brusoni et al_2001_ASQ_Knowledge Specialization Organizational Coupling and the Boundaries of the Firm.png'
$ws.Range("H31").Value = 'bar chart'
$ws.Range("I31").Value = $true

$ws.Range("A32").Value = 'Elsbach & Bhattacharya_2001_OrgSci_Defining Who you Are by what you''re not.pdf'
$ws.Range("B32").Value = 'F1_P14_Elsbach & Bhattacharya_2001_OrgSci_Defining Who you Are by what you''re not.png'
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 2001
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 'final_figures/2001/F1_P14_Elsbach & Bhattacharya_2001_OrgSci_Defining Who you Are by what you''re not.png'
$ws.Range("H32").Value = 'conceptual diagram'
$ws.Range("I32").Value = $true

$ws.Range("A33").Value = 'Ely & Thomas_2001_ASQ_Cultural Diversity at Work.pdf'
$ws.Range("B33").Value = 'F1_P7_Ely & Thomas_2001_ASQ_Cultural Diversity at Work.png'
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 2001
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 'final_figures/2001/F1_P7_Ely & Thomas_2001_ASQ_Cultural Diversity at Work.png'
$ws.Range("H33").Value = 'conceptual diagram'
$ws.Range("I33").Value = $true

$ws.Range("A34").Value = 'Hoffman & Ocasio_2001_OrgSci_Not All Events are Attended Equally.pdf'
$ws.Range("B34").Value = 'F1_P16_Hoffman & Ocasio_2001_OrgSci_Not All Events are Attended Equally.png'
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 2001
$ws.Range("E34").Value = 17
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 'final_figures/2001/F1_P16_Hoffman & Ocasio_2001_OrgSci_Not All Events are Attended Equally.png'
$ws.Range("H34").Value = 'process diagram'
$ws.Range("I34").Value = $true

$ws.Range("A35").Value = 'Mitzberg_2001_OrgSci_Managing Exceptionally.pdf'
$ws.Range("B35").Value = 'F1_P2_Mitzberg_2001_OrgSci_Managing Exceptionally.png'
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 2001
$ws.Range("E35").Value = 3
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 'final_figures/2001/F1_P2_Mitzberg_2001_OrgSci_Managing Exceptionally.png'
$ws.Range("H35").Value = 'conceptual diagram'
$ws.Range("I35").Value = $true

$ws.Range("A36").Value = 'Pitcher & Smith_2001_OrgSci_Top Management Team Heterogeneity_Quant.pdf'
$ws.Range("B36").Value = 'F1_P9_Pitcher & Smith_2001_OrgSci_Top Management Team Heterogeneity_Quant.png'
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 2001
$ws.Range("E36").Value = 10
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 'final_figures/2001/F1_P9_Pitcher & Smith_2001_OrgSci_Top Management Team Heterogeneity_Quant.png'
$ws.Range("H36").Value = 'line graph'
$ws.Range("I36").Value = $true
